{"js": "// Replace each two-digit multiplication expression with its updated\n// value. Every \"old\" value in the table is unique, so an exact,\n// non-wildcard search safely finds exactly one run to update.\nconst pairs = [\n    [\"13\u00d730=390\", \"40\u00d737=1480\"],\n    [\"80\u00d791=7280\", \"88\u00d742=3696\"],\n    [\"68\u00d737=2516\", \"97\u00d788=8536\"],\n    [\"54\u00d777=4158\", \"89\u00d728=2492\"],\n    [\"30\u00d751=1530\", \"88\u00d791=8008\"],\n    [\"78\u00d755=4290\", \"42\u00d795=3990\"],\n    [\"11\u00d741=451\", \"48\u00d789=4272\"],\n    [\"18\u00d737=666\", \"70\u00d766=4620\"],\n    [\"40\u00d721=840\", \"18\u00d795=1710\"],\n    [\"77\u00d741=3157\", \"31\u00d778=2418\"],\n    [\"78\u00d723=1794\", \"52\u00d732=1664\"],\n    [\"29\u00d796=2784\", \"52\u00d719=988\"],\n    [\"69\u00d760=4140\", \"28\u00d792=2576\"],\n    [\"46\u00d753=2438\", \"49\u00d758=2842\"],\n    [\"52\u00d787=4524\", \"33\u00d774=2442\"],\n    [\"95\u00d746=4370\", \"37\u00d763=2331\"],\n    [\"93\u00d783=7719\", \"46\u00d781=3726\"],\n    [\"38\u00d766=2508\", \"52\u00d782=4264\"],\n    [\"87\u00d785=7395\", \"38\u00d795=3610\"],\n    [\"39\u00d722=858\", \"88\u00d744=3872\"],\n    [\"99\u00d759=5841\", \"94\u00d741=3854\"],\n    [\"82\u00d718=1476\", \"27\u00d714=378\"],\n    [\"76\u00d786=6536\", \"66\u00d722=1452\"],\n    [\"26\u00d744=1144\", \"29\u00d771=2059\"],\n    [\"66\u00d747=3102\", \"74\u00d713=962\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression with its updated\n# value. Every \"old\" value in the table is unique, so an exact\n# Find/Replace (no wildcards) safely targets exactly one occurrence.\n$pairs = @(\n    @{ Find = '13\u00d730=390'; Replace = '40\u00d737=1480' },\n    @{ Find = '80\u00d791=7280'; Replace = '88\u00d742=3696' },\n    @{ Find = '68\u00d737=2516'; Replace = '97\u00d788=8536' },\n    @{ Find = '54\u00d777=4158'; Replace = '89\u00d728=2492' },\n    @{ Find = '30\u00d751=1530'; Replace = '88\u00d791=8008' },\n    @{ Find = '78\u00d755=4290'; Replace = '42\u00d795=3990' },\n    @{ Find = '11\u00d741=451'; Replace = '48\u00d789=4272' },\n    @{ Find = '18\u00d737=666'; Replace = '70\u00d766=4620' },\n    @{ Find = '40\u00d721=840'; Replace = '18\u00d795=1710' },\n    @{ Find = '77\u00d741=3157'; Replace = '31\u00d778=2418' },\n    @{ Find = '78\u00d723=1794'; Replace = '52\u00d732=1664' },\n    @{ Find = '29\u00d796=2784'; Replace = '52\u00d719=988' },\n    @{ Find = '69\u00d760=4140'; Replace = '28\u00d792=2576' },\n    @{ Find = '46\u00d753=2438'; Replace = '49\u00d758=2842' },\n    @{ Find = '52\u00d787=4524'; Replace = '33\u00d774=2442' },\n    @{ Find = '95\u00d746=4370'; Replace = '37\u00d763=2331' },\n    @{ Find = '93\u00d783=7719'; Replace = '46\u00d781=3726' },\n    @{ Find = '38\u00d766=2508'; Replace = '52\u00d782=4264' },\n    @{ Find = '87\u00d785=7395'; Replace = '38\u00d795=3610' },\n    @{ Find = '39\u00d722=858'; Replace = '88\u00d744=3872' },\n    @{ Find = '99\u00d759=5841'; Replace = '94\u00d741=3854' },\n    @{ Find = '82\u00d718=1476'; Replace = '27\u00d714=378' },\n    @{ Find = '76\u00d786=6536'; Replace = '66\u00d722=1452' },\n    @{ Find = '26\u00d744=1144'; Replace = '29\u00d771=2059' },\n    @{ Find = '66\u00d747=3102'; Replace = '74\u00d713=962' }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
